$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New course-taken records being appended for newly enrolled students.
$newRows = @(
    @{ A = "01-7877196"; B = "340x/SP16"; C = 1;    D = $null; E = $null },
    @{ A = "01-7877196"; B = "551x/SU17"; C = 1;    D = $null; E = $null },
    @{ A = "34-3115004"; B = "340x/SP16"; C = 0;    D = 1;    E = 1 },
    @{ A = "34-3115004"; B = "551x/SU17"; C = 1;    D = 1;    E = 0 },
    @{ A = "82-9118928"; B = "551x/SU17"; C = 0;    D = 1;    E = 1 }
)

$startRow = 200
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($row, 1).Value = $data.A
    $ws.Cells.Item($row, 2).Value = $data.B
    $ws.Cells.Item($row, 3).Value = $data.C

    if ($data.D -ne $null) {
        $ws.Cells.Item($row, 4).Value = $data.D
    }
    if ($data.E -ne $null) {
        $ws.Cells.Item($row, 5).Value = $data.E
    }
}

$ws.Range("D201").Select()
